# MCH242.xlsx update
# Adds a new data row (row 2) underneath the existing header row, giving
# the values for: identifier, levelOfDescription, extentAndMedium, notes.
# (columns B/C/D/H are left blank, matching the source record.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data -------------------------------------------------
$ws.Range("A2").Value = "MCH242"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"

# --- Formatting: match the font used for the rest of the sheet's data
# (10pt Calibri, automatic/text1 theme color) across the whole new row,
# including the blank cells so the row reads as one consistent record.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.ThemeColor = 1

# --- Keep the header row frozen and move the active selection down to
# the newly added row.
$ws.Range("A2:J2").Select()
$excel.ActiveWindow.FreezePanes = $true
